# The "imc" sheet classifies BMI (IMC) values in column D using a chain of
# nested IFs. The "Obesidad 2" threshold was mistyped as 29.9 (a duplicate
# of the "sobrepeso" threshold) instead of 39.9, so values between 34.9 and
# 39.9 were incorrectly skipping straight to "Obesidad 3". Fix the formula
# for every data row (2-100) on the "imc" sheet and leave the selection on D2,
# matching what Excel does after the user edits D2 and fills the correction
# down the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("imc")

for ($r = 2; $r -le 100; $r++) {
    $formula = '=IF(C' + $r + '<16,"Desnutrición severa",' +
        'IF(C' + $r + '<18.4,"Desnutrición moderada",' +
        'IF(C' + $r + '<22,"Bajo Peso",' +
        'IF(C' + $r + '<24.9,"Normal",' +
        'IF(C' + $r + '<29.9,"sobrepeso",' +
        'IF(C' + $r + '<34.9,"Obesidad 1",' +
        'IF(C' + $r + '<39.9,"Obesidad 2","Obesidad 3")))))))'
    $ws.Range("D$r").Formula = $formula
}

$ws.Range("D2").Select()
